$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Permutation: for target row index i (0-based, row = i+2), the value
# comes from original source row $perm[i]. Derived from the diff: the
# sheet's data rows (2..123) were reordered (re-sorted), carrying the
# columns Fecha(D), Volumen(J), Precio minimo(K), Precio maximo(L),
# Precio promedio ponderado(M), Origen(O) and Precio $/Kg(P) with them.
$perm = @(69,2,120,104,20,21,62,74,114,43,77,48,6,81,17,33,59,45,28,98,103,5,101,123,44,46,24,61,99,92,9,27,82,87,73,31,39,41,54,79,49,67,53,119,57,25,100,30,86,68,13,71,97,4,96,50,116,38,109,23,15,85,110,91,108,117,11,34,76,107,112,35,83,32,113,122,70,84,111,22,29,8,102,16,115,19,52,7,90,89,118,3,88,26,14,63,66,94,95,65,106,37,72,105,64,75,121,93,60,40,80,56,58,51,12,36,18,10,78,47,55,42)

$firstRow = 2
$lastRow = 123

# Snapshot the current ("before") values for the columns that move,
# keyed by row number, before any writes happen.
$snapD = @{}
$snapJ = @{}
$snapK = @{}
$snapL = @{}
$snapM = @{}
$snapO = @{}
$snapP = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, 4).Value2
    $snapJ[$r] = $ws.Cells.Item($r, 10).Value2
    $snapK[$r] = $ws.Cells.Item($r, 11).Value2
    $snapL[$r] = $ws.Cells.Item($r, 12).Value2
    $snapM[$r] = $ws.Cells.Item($r, 13).Value2
    $snapO[$r] = $ws.Cells.Item($r, 15).Value2
    $snapP[$r] = $ws.Cells.Item($r, 16).Value2
}

# Apply the permutation: row (i+2) gets the snapshot values that
# originally belonged to row $perm[i].
for ($i = 0; $i -lt $perm.Length; $i++) {
    $targetRow = $i + 2
    $srcRow = $perm[$i]

    $ws.Cells.Item($targetRow, 4).Value2 = $snapD[$srcRow]
    $ws.Cells.Item($targetRow, 10).Value2 = $snapJ[$srcRow]
    $ws.Cells.Item($targetRow, 11).Value2 = $snapK[$srcRow]
    $ws.Cells.Item($targetRow, 12).Value2 = $snapL[$srcRow]
    $ws.Cells.Item($targetRow, 13).Value2 = $snapM[$srcRow]
    $ws.Cells.Item($targetRow, 15).Value2 = $snapO[$srcRow]
    $ws.Cells.Item($targetRow, 16).Value2 = $snapP[$srcRow]
}

Write-Output "Permutation applied to rows $firstRow..$lastRow"
